$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")
$ws.Activate()

# ---------------------------------------------------------------------------
# 1. Add the new "Revenue" block in columns AD:AK.
#    (Text written first so the shared-string table keeps the same order as
#    the reference workbook: "Revenue" before "Evolution Gaming"/"EVO ST".)
# ---------------------------------------------------------------------------
# Row 2: merged "Revenue" header, styled like the existing Q2/W2 headers.
$ws.Range("Q2").Copy()
$ws.Range("AD2:AK2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("AD2").Value = "Revenue"
$ws.Range("AD2:AK2").Merge()

# ---------------------------------------------------------------------------
# 2. Insert a new row for "Evolution Gaming" right after Las Vegas Sands
#    (this shifts the old rows 5-17 down to 6-18, keeping their formatting).
# ---------------------------------------------------------------------------
$ws.Rows("5:5").Insert()

$ws.Range("C5").Value = "Evolution Gaming"
$ws.Range("D5").Value = "EVO ST"

# Rebuild the running index column (B) so formulas point at the row above
# them again (row insert does not auto-adjust the relative reference).
$ws.Range("B5").Formula = "=+B4+1"
for ($r = 6; $r -le 18; $r++) {
  $prev = $r - 1
  $ws.Range("B$r").Formula = "=+B$prev+1"
}

# Row 3: year headers 2019-2026 (first three are literals, the rest formulas).
$ws.Range("AD3").Value = 2019
$ws.Range("AE3").Value = 2020
$ws.Range("AF3").Value = 2021
$ws.Range("AD3:AF3").HorizontalAlignment = -4108

$ws.Range("AG3").Formula = "=+AF3+1"
$ws.Range("AH3:AK3").Formula = "=+AG3+1"

# Revenue figures (2019-2021) land on row 6 (the Flutter row after the insert).
$ws.Range("AD6").Value = 2140
$ws.Range("AE6").Value = 4414
$ws.Range("AF6").Value = 6036

# ---------------------------------------------------------------------------
# 3. Refresh the frozen panes / selection to include the new column D data.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("E4").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F5").Select()
